$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.988.99'
$ws.Range("E2").Value = '  +0.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.373.50'
$ws.Range("E3").Value = '  -3.07%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.57'
$ws.Range("E5").Value = '  -0.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.76'
$ws.Range("E6").Value = '  +3.61%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.627'
$ws.Range("E7").Value = '  +5.16%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.370.22'
$ws.Range("E9").Value = '  -3.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.131'
$ws.Range("E10").Value = '  -0.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.92'
$ws.Range("E11").Value = '  +0.98%  '
$ws.Range("E12").Value = '  +0.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.958.71'
$ws.Range("E13").Value = '  -3.15%  '
$ws.Range("E14").Value = '  +1.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.95'
$ws.Range("E15").Value = '  -3.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.934.68'
$ws.Range("E17").Value = '  -0.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.389.92'
$ws.Range("E18").Value = '  -2.62%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.81'
$ws.Range("E19").Value = '  -2.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.61'
$ws.Range("E20").Value = '  -2.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '365.64'
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.52'
$ws.Range("E22").Value = '  -3.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.52'
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("E24").Value = '  -0.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.529'
$ws.Range("E25").Value = '  -1.07%  '
$ws.Range("E26").Value = '  +0.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.74'
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("E28").Value = '  +0.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("E30").Value = '  -0.33%  '
$ws.Range("E31").Value = '  -0.84%  '
$ws.Range("E32").Value = '  -4.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.97'
$ws.Range("E34").Value = '  -2.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.24'
$ws.Range("E35").Value = '  -4.02%  '
$ws.Range("E36").Value = '  -1.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '160.97'
$ws.Range("E37").Value = '  +0.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.856'
$ws.Range("E38").Value = '  -3.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '27.02'
$ws.Range("E39").Value = '  -8.52%  '
$ws.Range("E40").Value = '  +0.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.60'
$ws.Range("E41").Value = '  +1.33%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.686.32'
$ws.Range("E42").Value = '  -4.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.32'
$ws.Range("E43").Value = '  -0.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.33'
$ws.Range("E44").Value = '  -3.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0674'
$ws.Range("E45").Value = '  -1.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '338.62'
$ws.Range("E46").Value = '  +10.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '39.87'
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.33'
$ws.Range("E48").Value = '  +0.81%  '
$ws.Range("E49").Value = '  -2.00%  '
$ws.Range("E50").Value = '  +3.11%  '
$ws.Range("B51").Value = 'ONDO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.983'
$ws.Range("E51").Value = '  +0.44%  '
